# Update the cryptocurrency "Price" (column D) values to reflect the
# latest pull of symbol data, as published by the GitHub Actions job.
# The sheet stores these prices as plain text (not numbers), so we force
# a Text number format while writing to avoid Excel auto-converting the
# numeric-looking strings into real numbers, then restore each cell's
# original style so no visual formatting is changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "268.36"
    3  = "22.93"
    4  = "6.354"
    5  = "0.06228"
    6  = "3.615"
    7  = "6.695"
    9  = "0.8304"
    10 = "0.01365"
    11 = "0.1616"
    12 = "0.08279"
    13 = "0.03411"
    14 = "0.03175"
    15 = "0.09310"
    16 = "3.930"
    17 = "0.001714"
    18 = "0.04853"
    19 = "0.006301"
    20 = "0.005401"
    21 = "0.001092"
    23 = "3.769"
    24 = "2.366"
    25 = "0.3351"
    26 = "0.1244"
    27 = "0.0002688"
    40 = "0.04668"
    41 = "0.006896"
    42 = "0.1161"
    43 = "0.003468"
    44 = "0.01237"
    45 = "0.00006266"
    46 = "0.00000000751"
    47 = "0.7905"
    48 = "0.1508"
    49 = "0.00002104"
    50 = "0.01242"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Range("D$row")
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = $origStyle
}
